$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the header text for the "Max Incoming/Outgoing Interactions" columns
# (C1 <-> E1), reflecting the new cluster order in heatmaps.
$ws.Range("C1").Value = "Max Outgoing Interactions"
$ws.Range("E1").Value = "Max Incoming Interactions"
